$d = $word.ActiveDocument

# The target paragraph currently reads:
#   "Observer sollen jedoch zusätzlich zu "  (one run)
# and must become four separate runs (same rPr throughout):
#   "Observer" | "*innen" | " sollen jedoch zusätzlich zu " | "den "
# immediately followed by the untouched "Observations" run (still wrapped
# in its proofErr spell-check markers).

# Step 1: replace the whole original phrase with the fully-combined text.
# This lands as a single run (no split yet), narrowed exactly to the new text.
$rng = $d.Content
$found = $rng.Find.Execute(
    "Observer sollen jedoch zusätzlich zu ",
    $true, $false, $false, $false, $false, $true, 0, $false,
    "Observer*innen sollen jedoch zusätzlich zu den ", 1)

if (-not $found) {
    throw "Could not find target phrase to replace."
}

$segStart = $rng.Start
$segEnd = $rng.End

# Step 2: split that merged run into the four desired pieces by toggling a
# character property on each sub-range (on, then immediately back off) —
# this forces Word to break run boundaries exactly at the sub-range edges
# without altering any visible formatting.

$r1 = $d.Range($segStart, $segEnd)
$r1.Find.Execute("*innen", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$r1.Bold = 1
$r1.Bold = 0

$r2 = $d.Range($segStart, $segEnd)
$r2.Find.Execute(" sollen jedoch zusätzlich zu ", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$r2.Bold = 1
$r2.Bold = 0

$r3 = $d.Range($segStart, $segEnd)
$r3.Find.Execute("den ", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$r3.Bold = 1
$r3.Bold = 0
